$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 15: correct typo "Realisation and Desing" -> "Realisation and Design" ---
$ws.Range("E15").Value = "Realisation and Design"

# --- Row 16 ---
$ws.Range("A16").Value = 40247
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B16").Value = 0.39583333333333331
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C16").Value = 0.64583333333333337
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 6
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "Realisation and Design"
$ws.Range("F16").Value = "Design multiple kinds of surfaces, enemy implementation"

# --- Row 17 ---
$ws.Range("A17").Value = 40248
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B17").Value = 0.4375
$ws.Range("B15").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C17").Value = 0.64583333333333337
$ws.Range("C15").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D17").Value = 5
$ws.Range("D15").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "Realisation and Design"
$ws.Range("F17").Value = "Enemy factory"

# --- Row 18 ---
$ws.Range("A18").Value = 40252
$ws.Range("A15").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B18").Value = 0.45833333333333331
$ws.Range("B15").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C18").Value = 0.60416666666666663
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
# D18 hours stored as TEXT "4.5" (quirk matching original file, e.g. D8/D12)
$ws.Range("D18").Formula = '="4.5"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("D15").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "Design"
$ws.Range("F18").Value = "Tiles for the snow level"

# --- Row 19 ---
$ws.Range("A19").Value = 40253
$ws.Range("A15").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B19").Value = 0.42708333333333331
$ws.Range("B15").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C19").Value = 0.625
$ws.Range("C15").Copy()
$ws.Range("C19").PasteSpecial(-4122)
# D19 hours stored as TEXT "4.75"
$ws.Range("D19").Formula = '="4.75"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("D15").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "Realisation"
$ws.Range("F19").Value = "Snowlevel implementation, all tiles and some surfaces"

# --- Rows 20-70: only column D has an empty cell with the right-aligned number style ---
$ws.Range("D15").Copy()
$ws.Range("D20:D70").PasteSpecial(-4122)
$ws.Range("D20:D70").ClearContents()

# --- Update the sheet view (scrolled position + active selection) ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F21").Select()

Write-Host "done"
